# Updated cryptos list on Tue Oct 17 07:17:43 UTC 2023 with GitHub Actions
# Refreshes Price (D) and Volume/1h (E) text cells for the coinranking rows
# that moved since the previous snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.255.70"
$ws.Range("E2").Value = "  +2.61%  "
$ws.Range("D3").Value = "1.586.50"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  +1.12%  "
$ws.Range("D5").Value = "'213.45"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D8").Value = "'24.00"
$ws.Range("E8").Value = "  +7.04%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").Value = "1.812.83"
$ws.Range("E12").Value = "  +1.48%  "
$ws.Range("D13").Value = "1.585.83"
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("E14").Value = "  +1.90%  "
$ws.Range("D15").Value = "'3.75"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "28.282.24"
$ws.Range("E16").Value = "  +2.79%  "
$ws.Range("D17").Value = "'63.17"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").Value = "'227.58"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("D20").Value = "'7.47"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("E21").Value = "  +1.03%  "
$ws.Range("E22").Value = "  -1.54%  "
$ws.Range("D23").Value = "'9.32"
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").Value = "'151.94"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").Value = "'15.14"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").Value = "'3.24"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D34").Value = "1.399.37"
$ws.Range("E34").Value = "  -4.35%  "
$ws.Range("E35").Value = "  -1.27%  "
$ws.Range("E36").Value = "  -7.74%  "
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("E39").Value = "  +8.87%  "
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("E41").Value = "  -0.58%  "
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").Value = "'1.89"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").Value = "'5.60"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").Value = "'64.25"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "1.723.19"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("D49").Value = "'86.59"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("E50").Value = "  +2.83%  "
$ws.Range("D51").Value = "'0.0519"
$ws.Range("E51").Value = "  -0.71%  "
